$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 61104
$ws.Range("E2").Value = 1203494574653
$ws.Range("F2").Value = 12363680338
$ws.Range("G2").Value = 0.56187

# Row 3
$ws.Range("D3").Value = 2926.93
$ws.Range("E3").Value = 351603917193
$ws.Range("F3").Value = 5634080405
$ws.Range("G3").Value = 0.7301800000000001

# Row 4
$ws.Range("D4").Value = 0.99978
$ws.Range("E4").Value = 110823059182
$ws.Range("F4").Value = 20087751557
$ws.Range("G4").Value = 0.03389

# Row 5
$ws.Range("D5").Value = 593.3200000000001
$ws.Range("E5").Value = 91270299632
$ws.Range("F5").Value = 515291294
$ws.Range("G5").Value = 0.60704

# Row 6
$ws.Range("D6").Value = 145.21
$ws.Range("E6").Value = 65128845258
$ws.Range("F6").Value = 1118899080
$ws.Range("G6").Value = 0.40865

# Row 7
$ws.Range("D7").Value = 1
$ws.Range("E7").Value = 33130304606
$ws.Range("F7").Value = 2732154575
$ws.Range("G7").Value = 0.00546

# Row 8
$ws.Range("B8").Value = "XRP"
$ws.Range("C8").Value = "XRP"
$ws.Range("D8").Value = 0.504237
$ws.Range("E8").Value = 27910099562
$ws.Range("F8").Value = 346462825
$ws.Range("G8").Value = 0.20126

# Row 9
$ws.Range("B9").Value = "STETH"
$ws.Range("C9").Value = "Lido Staked Ether"
$ws.Range("D9").Value = 2925.7
$ws.Range("E9").Value = 27401981511
$ws.Range("F9").Value = 27041451
$ws.Range("G9").Value = 0.7047099999999999

# Row 10
$ws.Range("B10").Value = "TON"
$ws.Range("C10").Value = "Toncoin"
$ws.Range("D10").Value = 7.03
$ws.Range("E10").Value = 24371890771
$ws.Range("F10").Value = 280741207
$ws.Range("G10").Value = 4.90691

# Row 11
$ws.Range("B11").Value = "DOGE"
$ws.Range("C11").Value = "Dogecoin"
$ws.Range("D11").Value = 0.142743
$ws.Range("E11").Value = 20626235301
$ws.Range("F11").Value = 607572981
$ws.Range("G11").Value = -0.64562

# Row 12
$ws.Range("D12").Value = 0.439853
$ws.Range("E12").Value = 15559954945
$ws.Range("F12").Value = 161066423
$ws.Range("G12").Value = -0.86084

# Row 13
$ws.Range("D13").Value = 0.00002251
$ws.Range("E13").Value = 13282895467
$ws.Range("F13").Value = 187059494
$ws.Range("G13").Value = -0.29554

# Row 14
$ws.Range("D14").Value = 33.7
$ws.Range("E14").Value = 12881563701
$ws.Range("F14").Value = 205379624
$ws.Range("G14").Value = 0.83646

# Row 15
$ws.Range("D15").Value = 0.126671
$ws.Range("E15").Value = 11080231569
$ws.Range("F15").Value = 222040447
$ws.Range("G15").Value = 0.04657

# Row 16
$ws.Range("D16").Value = 61145
$ws.Range("E16").Value = 9491996464
$ws.Range("F16").Value = 81669094
$ws.Range("G16").Value = 0.65043

# Row 17
$ws.Range("D17").Value = 6.73
$ws.Range("E17").Value = 9185077286
$ws.Range("F17").Value = 99309899
$ws.Range("G17").Value = 0.37324

# Row 18
$ws.Range("D18").Value = 435.06
$ws.Range("E18").Value = 8576584447
$ws.Range("F18").Value = 126456076
$ws.Range("G18").Value = 1.31722

# Row 19
$ws.Range("D19").Value = 13.44
$ws.Range("E19").Value = 7900225468
$ws.Range("F19").Value = 188113916
$ws.Range("G19").Value = -0.57428

# Row 20
$ws.Range("D20").Value = 6.95
$ws.Range("E20").Value = 7453224300
$ws.Range("F20").Value = 252845305
$ws.Range("G20").Value = -1.70197

# Row 21
$ws.Range("D21").Value = 0.678713
$ws.Range("E21").Value = 6306894989
$ws.Range("F21").Value = 137243398
$ws.Range("G21").Value = -0.48911

# Row 22
$ws.Range("D22").Value = 81.48
$ws.Range("E22").Value = 6072821910
$ws.Range("F22").Value = 202169810
$ws.Range("G22").Value = -0.57816

# Row 23
$ws.Range("B23").Value = "FET"
$ws.Range("C23").Value = "Fetch.ai"
$ws.Range("D23").Value = 2.2
$ws.Range("E23").Value = 5558569059
$ws.Range("F23").Value = 125335431
$ws.Range("G23").Value = -0.78065

# Row 24
$ws.Range("B24").Value = "ICP"
$ws.Range("C24").Value = "Internet Computer"
$ws.Range("D24").Value = 11.83
$ws.Range("E24").Value = 5494792187
$ws.Range("F24").Value = 59699771
$ws.Range("G24").Value = -0.92692

# Row 25
$ws.Range("B25").Value = "DAI"
$ws.Range("C25").Value = "Dai"
$ws.Range("D25").Value = 0.999709
$ws.Range("E25").Value = 5473476264
$ws.Range("F25").Value = 634862514
$ws.Range("G25").Value = -0.01465

# Row 26
$ws.Range("B26").Value = "LEO"
$ws.Range("C26").Value = "LEO Token"
$ws.Range("D26").Value = 5.9
$ws.Range("E26").Value = 5459889666
$ws.Range("F26").Value = 1568276
$ws.Range("G26").Value = 0.31157

# Row 27
$ws.Range("B27").Value = "UNI"
$ws.Range("C27").Value = "Uniswap"
$ws.Range("D27").Value = 7.11
$ws.Range("E27").Value = 5363060020
$ws.Range("F27").Value = 86439897
$ws.Range("G27").Value = 0.52007

# Row 28
$ws.Range("B28").Value = "RNDR"
$ws.Range("C28").Value = "Render"
$ws.Range("D28").Value = 11.03
$ws.Range("E28").Value = 4296394894
$ws.Range("F28").Value = 198928600
$ws.Range("G28").Value = 1.92776

# Row 29
$ws.Range("D29").Value = 0.110663
$ws.Range("E29").Value = 3949696675
$ws.Range("F29").Value = 53147580
$ws.Range("G29").Value = 3.23655

# Row 30
$ws.Range("B30").Value = "ETC"
$ws.Range("C30").Value = "Ethereum Classic"
$ws.Range("D30").Value = 26.69
$ws.Range("E30").Value = 3925985861
$ws.Range("F30").Value = 90448714
$ws.Range("G30").Value = 0.83439

# Row 31
$ws.Range("D31").Value = 1
$ws.Range("E31").Value = 3815889177
$ws.Range("F31").Value = 3224194212
$ws.Range("G31").Value = -0.02008

# Row 32
$ws.Range("B32").Value = "PEPE"
$ws.Range("C32").Value = "Pepe"
$ws.Range("D32").Value = 0.00000867
$ws.Range("E32").Value = 3649848628
$ws.Range("F32").Value = 429447213
$ws.Range("G32").Value = 2.22099

# Row 33
$ws.Range("B33").Value = "APT"
$ws.Range("C33").Value = "Aptos"
$ws.Range("D33").Value = 8.41
$ws.Range("E33").Value = 3611222310
$ws.Range("F33").Value = 48373512
$ws.Range("G33").Value = 0.15684

# Row 34
$ws.Range("B34").Value = "ATOM"
$ws.Range("C34").Value = "Cosmos Hub"
$ws.Range("D34").Value = 8.58
$ws.Range("E34").Value = 3355493362
$ws.Range("F34").Value = 99896835
$ws.Range("G34").Value = -0.05023

# Row 35
$ws.Range("B35").Value = "CRO"
$ws.Range("C35").Value = "Cronos"
$ws.Range("D35").Value = 0.124056
$ws.Range("E35").Value = 3315125367
$ws.Range("F35").Value = 8115194
$ws.Range("G35").Value = -0.24031

# Row 36
$ws.Range("D36").Value = 1.01
$ws.Range("E36").Value = 3298484571
$ws.Range("F36").Value = 45519080
$ws.Range("G36").Value = 0.15957

# Row 37
$ws.Range("B37").Value = "IMX"
$ws.Range("C37").Value = "Immutable"
$ws.Range("D37").Value = 2.25
$ws.Range("E37").Value = 3274684827
$ws.Range("F37").Value = 58091866
$ws.Range("G37").Value = 0.26236

# Row 38
$ws.Range("B38").Value = "WEETH"
$ws.Range("C38").Value = "Wrapped eETH"
$ws.Range("D38").Value = 3038.88
$ws.Range("E38").Value = 3127787678
$ws.Range("F38").Value = 20156069
$ws.Range("G38").Value = 0.787

# Row 39
$ws.Range("D39").Value = 5.64
$ws.Range("E39").Value = 3112314550
$ws.Range("F39").Value = 76316825
$ws.Range("G39").Value = 0.9496

# Row 40
$ws.Range("B40").Value = "XLM"
$ws.Range("C40").Value = "Stellar"
$ws.Range("D40").Value = 0.105414
$ws.Range("E40").Value = 3051115291
$ws.Range("F40").Value = 30043559
$ws.Range("G40").Value = -0.80449

# Row 41
$ws.Range("B41").Value = "OKB"
$ws.Range("C41").Value = "OKB"
$ws.Range("D41").Value = 49.94
$ws.Range("E41").Value = 2995043075
$ws.Range("F41").Value = 2712624
$ws.Range("G41").Value = 1.08076

# Row 42
$ws.Range("B42").Value = "WIF"
$ws.Range("C42").Value = "dogwifhat"
$ws.Range("D42").Value = 2.98
$ws.Range("E42").Value = 2988678110
$ws.Range("F42").Value = 291328777
$ws.Range("G42").Value = -0.85681

# Row 43
$ws.Range("B43").Value = "KAS"
$ws.Range("C43").Value = "Kaspa"
$ws.Range("D43").Value = 0.124618
$ws.Range("E43").Value = 2934998933
$ws.Range("F43").Value = 21796106
$ws.Range("G43").Value = 0.08401

# Row 44
$ws.Range("B44").Value = "STX"
$ws.Range("C44").Value = "Stacks"
$ws.Range("D44").Value = 1.99
$ws.Range("E44").Value = 2908293840
$ws.Range("F44").Value = 26951589
$ws.Range("G44").Value = -0.67747

# Row 45
$ws.Range("B45").Value = "EZETH"
$ws.Range("C45").Value = "Renzo Restaked ETH"
$ws.Range("D45").Value = 2880.77
$ws.Range("E45").Value = 2870481798
$ws.Range("F45").Value = 26251946
$ws.Range("G45").Value = 0.76324

# Row 46
$ws.Range("B46").Value = "AR"
$ws.Range("C46").Value = "Arweave"
$ws.Range("D46").Value = 42.08
$ws.Range("E46").Value = 2749163283
$ws.Range("F46").Value = 98722764
$ws.Range("G46").Value = 3.42502

# Row 47
$ws.Range("B47").Value = "GRT"
$ws.Range("C47").Value = "The Graph"
$ws.Range("D47").Value = 0.284524
$ws.Range("E47").Value = 2702710102
$ws.Range("F47").Value = 77186920
$ws.Range("G47").Value = -3.46272

# Row 48
$ws.Range("B48").Value = "OP"
$ws.Range("C48").Value = "Optimism"
$ws.Range("D48").Value = 2.56
$ws.Range("E48").Value = 2676668928
$ws.Range("F48").Value = 91814297
$ws.Range("G48").Value = 0.22292

# Row 49
$ws.Range("B49").Value = "ARB"
$ws.Range("C49").Value = "Arbitrum"
$ws.Range("D49").Value = 1.002
$ws.Range("E49").Value = 2661095585
$ws.Range("F49").Value = 130334444
$ws.Range("G49").Value = 0.50502

# Row 50
$ws.Range("B50").Value = "TAO"
$ws.Range("C50").Value = "Bittensor"
$ws.Range("D50").Value = 374.54
$ws.Range("E50").Value = 2534911721
$ws.Range("F50").Value = 16741618
$ws.Range("G50").Value = 0.2283

# Row 51
$ws.Range("B51").Value = "VET"
$ws.Range("C51").Value = "VeChain"
$ws.Range("D51").Value = 0.03462744
$ws.Range("E51").Value = 2518990830
$ws.Range("F51").Value = 39645428
$ws.Range("G51").Value = -0.97274
